# The workbook ships with exactly one sheet: "Wahlpflichtmodul".
# Rename it to match the new Threagile test model.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Insecure IT Testmodell"

# Column B header changes from the old tag name to the new "pii" tag.
$ws.Cells.Item(1, 2).Value = "pii"

# New tag rows (column A) and whether the "pii" tag applies (column B = X).
$tags = @(
    @{ Name = "app";         Marked = $false },
    @{ Name = "to-database"; Marked = $false },
    @{ Name = "client";      Marked = $false },
    @{ Name = "to-webapp";   Marked = $false },
    @{ Name = "database";    Marked = $false },
    @{ Name = "webapp";      Marked = $false },
    @{ Name = "to-app";      Marked = $false },
    @{ Name = "user-data";   Marked = $true  },
    @{ Name = "dmz";         Marked = $false }
)

# Rows 2-9 already exist with the right formatting; row 10 is brand new, so
# clone the formatting from row 9 before writing into it.
$ws.Range("A9:C9").Copy() | Out-Null
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$row = 2
foreach ($tag in $tags) {
    $ws.Cells.Item($row, 1).Value = $tag.Name
    if ($tag.Marked) {
        $ws.Cells.Item($row, 2).Value = "X"
    } else {
        $ws.Cells.Item($row, 2).Value = $null
    }
    $row = $row + 1
}

# The matrix now only needs the "Element" and "pii" columns; drop column C.
$ws.Columns.Item(3).Delete()

# Note: the sheet's "first page" header (Page Setup > Header/Footer > Custom
# Header, "first page" variant) also embeds the sheet title, but Excel's
# object model only exposes Left/Center/RightHeader for the default (odd)
# header/footer - there is no VBA/COM property for the first-page header, so
# it cannot be updated from here.
